# Updated symbol list on Sat Feb 11 14:14:27 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume/hour snapshot rows (2-51) with the
# latest scraped values from coinranking.com. A handful of rows (6-17)
# also pick up reshuffled coin/link pairs as the source ranking moved.
#
# All cells in this sheet are plain text ("inlineStr") - not numbers or
# percentages - so each target cell is forced to Text format ("@") before
# the value is written. Otherwise Excel would auto-convert numeric-looking
# strings (prices, "13"/"14" hour values) and percent-looking strings
# ("0.73%") into real numbers, which would not match the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.73%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '14'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '5.56%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '14'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.120'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.56%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '14'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07637'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.46%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '14'

# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.619'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.76%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '14'

# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'BTSEToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '2.488'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.21%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '14'

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9058'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.97%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '14'

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1116'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '10.07%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '14'

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1805'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.91%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '14'

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09118'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.66%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '14'

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04252'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.00%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '14'

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1052'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.25%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '14'

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001258'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.42%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '14'

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005870'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.28%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '14'

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.345'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.49%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '14'

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.266'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.41%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '14'

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.32%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '14'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.630'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-5.32%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '14'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1359'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.79%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '14'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2736'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.89%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '14'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04042'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.07%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '14'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001265'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5.01%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '14'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004054'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.10%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '14'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001301'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.01%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '14'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003746'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '14'

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '14'

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '14'

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '14'

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '14'

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '14'

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '14'

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '14'

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '14'

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '14'

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '14'

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '14'

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02425'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '0.21%'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '14'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05255'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.64%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '14'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007780'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.61%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '14'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1303'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.99%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '14'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006618'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-6.95%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '14'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001951'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.02%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '14'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007559'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-9.84%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '14'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3085'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-7.66%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '14'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006790'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '5.89%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '14'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '14'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06204'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1,343.35%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '14'

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '39.92%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '14'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '14'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.01%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '14'
